$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.46
$ws.Range("G2").Value = 1.63
$ws.Range("H2").Value = 6.6
$ws.Range("J2").Value = 3.55
$ws.Range("K2").Value = 5.5
$ws.Range("L2").Value = 1.35
$ws.Range("S2").Value = 3.3
$ws.Range("T2").Value = 2.12
$ws.Range("U2").Value = 1.69
$ws.Range("W2").Value = 2.66
$ws.Range("F3").Value = 1.84
$ws.Range("G3").Value = 2.04
$ws.Range("I3").Value = 5.8
$ws.Range("J3").Value = 3.15
$ws.Range("K3").Value = 3.8
$ws.Range("L3").Value = 1.47
$ws.Range("P3").Value = 1.66
$ws.Range("Q3").Value = 2.2
$ws.Range("U3").Value = 1.82
$ws.Range("V3").Value = 1.21
$ws.Range("W3").Value = 1.96
$ws.Range("Y3").Value = 17.5
$ws.Range("Z3").Value = 46
$ws.Range("AC3").Value = 9.4
$ws.Range("F4").Value = 1.75
$ws.Range("G4").Value = 1.87
$ws.Range("H4").Value = 5.5
$ws.Range("I4").Value = 6.8
$ws.Range("J4").Value = 3.35
$ws.Range("K4").Value = 3.75
$ws.Range("N4").Value = 2.66
$ws.Range("P4").Value = 1.56
$ws.Range("V4").Value = 1.18
$ws.Range("W4").Value = 2.14
$ws.Range("AA4").Value = 1000
$ws.Range("H5").Value = 1.52
$ws.Range("I5").Value = 1.53
$ws.Range("J5").Value = 4.6
$ws.Range("K5").Value = 4.7
$ws.Range("N5").Value = 4.4
$ws.Range("O5").Value = 1.27
$ws.Range("P5").Value = 2.18
$ws.Range("Q5").Value = 1.82
$ws.Range("R5").Value = 1.46
$ws.Range("S5").Value = 3.05
$ws.Range("U5").Value = 1.99
$ws.Range("V5").Value = 2.88
$ws.Range("X5").Value = 17
$ws.Range("AC5").Value = 10
$ws.Range("AD5").Value = 9.800000000000001
$ws.Range("AJ5").Value = 240
$ws.Range("AL5").Value = 100
$ws.Range("AO5").Value = 7.4
$ws.Range("F6").Value = 1.25
$ws.Range("K6").Value = 8
$ws.Range("L6").Value = 1.22
$ws.Range("R6").Value = 1.97
$ws.Range("S6").Value = 2
$ws.Range("T6").Value = 1.88
$ws.Range("U6").Value = 2.02
$ws.Range("W6").Value = 4.8
$ws.Range("Y6").Value = 650
$ws.Range("AB6").Value = 14.5
$ws.Range("AE6").Value = 180
$ws.Range("AF6").Value = 10
$ws.Range("AJ6").Value = 10
$ws.Range("AN6").Value = 3.15
$ws.Range("AO6").Value = 140
$ws.Range("H7").Value = 2.3
$ws.Range("K7").Value = 3.7
$ws.Range("R7").Value = 1.51
$ws.Range("T7").Value = 1.62
$ws.Range("U7").Value = 2.52
$ws.Range("V7").Value = 1.76
$ws.Range("W7").Value = 1.41
$ws.Range("X7").Value = 17.5
$ws.Range("Y7").Value = 13
$ws.Range("Z7").Value = 16
$ws.Range("AA7").Value = 29
$ws.Range("AB7").Value = 16
$ws.Range("AD7").Value = 11
$ws.Range("AE7").Value = 19
$ws.Range("AF7").Value = 25
$ws.Range("AH7").Value = 15
$ws.Range("AN7").Value = 26
$ws.Range("AO7").Value = 14.5
$ws.Range("F8").Value = 1.46
$ws.Range("G8").Value = 1.47
$ws.Range("H8").Value = 8.199999999999999
$ws.Range("J8").Value = 5.1
$ws.Range("K8").Value = 5.2
$ws.Range("Q8").Value = 1.83
$ws.Range("S8").Value = 3.15
$ws.Range("T8").Value = 2.06
$ws.Range("W8").Value = 3.1
$ws.Range("X8").Value = 18.5
$ws.Range("AA8").Value = 300
$ws.Range("AD8").Value = 30
$ws.Range("AE8").Value = 130
$ws.Range("AG8").Value = 9.6
$ws.Range("AN8").Value = 7.2
$ws.Range("AO8").Value = 150
$ws.Range("H9").Value = 2.34
$ws.Range("I9").Value = 2.36
$ws.Range("J9").Value = 3.9
$ws.Range("K9").Value = 3.95
$ws.Range("P9").Value = 2.64
$ws.Range("V9").Value = 1.74
$ws.Range("W9").Value = 1.46
$ws.Range("X9").Value = 23
$ws.Range("AA9").Value = 32
$ws.Range("AC9").Value = 9.199999999999999
$ws.Range("AO9").Value = 11.5
$ws.Range("F10").Value = 2.28
$ws.Range("G10").Value = 2.32
$ws.Range("H10").Value = 3.2
$ws.Range("I10").Value = 3.25
$ws.Range("M10").Value = 1.04
$ws.Range("P10").Value = 2.66
$ws.Range("U10").Value = 2.78
$ws.Range("V10").Value = 1.44
$ws.Range("W10").Value = 1.76
$ws.Range("X10").Value = 23
$ws.Range("Y10").Value = 20
$ws.Range("AO10").Value = 18.5
$ws.Range("F11").Value = 2.18
$ws.Range("G11").Value = 2.2
$ws.Range("H11").Value = 3.55
$ws.Range("I11").Value = 3.65
$ws.Range("N11").Value = 4.8
$ws.Range("Q11").Value = 1.73
$ws.Range("V11").Value = 1.37
$ws.Range("W11").Value = 1.83
$ws.Range("AD11").Value = 14.5
$ws.Range("I12").Value = 21
$ws.Range("X12").Value = 900
$ws.Range("AB12").Value = 19.5
$ws.Range("AF12").Value = 12.5
$ws.Range("AG12").Value = 13.5
$ws.Range("AL12").Value = 30
$ws.Range("AO12").Value = 200
$ws.Range("I13").Value = 2.42
$ws.Range("P13").Value = 2.32
$ws.Range("T13").Value = 1.61
$ws.Range("U13").Value = 2.56
$ws.Range("V13").Value = 1.7
$ws.Range("W13").Value = 1.46
$ws.Range("AI13").Value = 32
$ws.Range("AO13").Value = 14.5
$ws.Range("F14").Value = 2.36
$ws.Range("G14").Value = 2.66
$ws.Range("H14").Value = 2.74
$ws.Range("I14").Value = 3.35
$ws.Range("J14").Value = 3.65
$ws.Range("L14").Value = 1.34
$ws.Range("M14").Value = 1.04
$ws.Range("N14").Value = 2.86
$ws.Range("O14").Value = 1.25
$ws.Range("P14").Value = 1.92
$ws.Range("Q14").Value = 1.58
$ws.Range("R14").Value = 1.39
$ws.Range("S14").Value = 2.28
$ws.Range("T14").Value = 1.04
$ws.Range("U14").Value = 1.9
$ws.Range("V14").Value = 1.44
$ws.Range("W14").Value = 1.6
$ws.Range("X14").Value = 1000
$ws.Range("Y14").Value = 20
$ws.Range("Z14").Value = 30
$ws.Range("AA14").Value = 1000
$ws.Range("AB14").Value = 1000
$ws.Range("AC14").Value = 1000
$ws.Range("AD14").Value = 19
$ws.Range("AE14").Value = 1000
$ws.Range("AF14").Value = 1000
$ws.Range("AG14").Value = 1000
$ws.Range("AH14").Value = 1000
$ws.Range("AI14").Value = 1000
$ws.Range("AJ14").Value = 1000
$ws.Range("AK14").Value = 36
$ws.Range("AL14").Value = 1000
$ws.Range("AM14").Value = 1000
$ws.Range("AN14").Value = 1000
$ws.Range("AO14").Value = 1000
